$wb = $excel.ActiveWorkbook

# Select the sheet "4x4 in 6x6 Deflated At 4" and update its selection
$wsDeflated4 = $wb.Worksheets.Item("4x4 in 6x6 Deflated At 4")
$wsDeflated4.Activate()
$wsDeflated4.Range("G4:J5").Select()

# Select the sheet "4x4 in 6x6" and update its selection, making it the active (last-active) tab
$ws6x6 = $wb.Worksheets.Item("4x4 in 6x6")
$ws6x6.Activate()
$ws6x6.Cells.Select()
$ws6x6.Range("I23").Select()
$ws6x6.Cells.Select()

# Force a recalculation so RAND() formulas on "rand" sheet produce new volatile values
$excel.CalculateFull()
